$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge the new row block (rows 18:23) per column, matching the target layout ---
$ws.Range("B18:B23").Merge()
$ws.Range("C18:C23").Merge()
$ws.Range("D18:D23").Merge()
$ws.Range("E18:E23").Merge()
$ws.Range("F18:F23").Merge()
$ws.Range("G18:G23").Merge()
$ws.Range("H18:H23").Merge()

# --- Alignment: B,C,D,F,G,H columns = center/center; E column = left/center + wrap (same as row 13-15 pattern) ---
$ws.Range("B18:B23").HorizontalAlignment = -4108
$ws.Range("B18:B23").VerticalAlignment = -4108
$ws.Range("B18:B23").WrapText = $false
$ws.Range("C18:C23").HorizontalAlignment = -4108
$ws.Range("C18:C23").VerticalAlignment = -4108
$ws.Range("C18:C23").WrapText = $false
$ws.Range("D18:D23").HorizontalAlignment = -4108
$ws.Range("D18:D23").VerticalAlignment = -4108
$ws.Range("D18:D23").WrapText = $false
$ws.Range("F18:F23").HorizontalAlignment = -4108
$ws.Range("F18:F23").VerticalAlignment = -4108
$ws.Range("F18:F23").WrapText = $false
$ws.Range("G18:G23").HorizontalAlignment = -4108
$ws.Range("G18:G23").VerticalAlignment = -4108
$ws.Range("G18:G23").WrapText = $false
$ws.Range("H18:H23").HorizontalAlignment = -4108
$ws.Range("H18:H23").VerticalAlignment = -4108
$ws.Range("H18:H23").WrapText = $false
$eRng = $ws.Range("E18:E23")
$eRng.HorizontalAlignment = -4131
$eRng.VerticalAlignment = -4108
$eRng.WrapText = $true

# --- Cell values (writes new shared strings "Chungbuk University - NCLab" then the long post) ---
$ws.Range("B18").Value = 5
$ws.Range("C18").Value = "Chungbuk University - NCLab"
$longText = @'
Chào mọi người,
Lab mình (NCLab- Chungbuk National University) đang tuyển thêm sinh viên PhD, Combined và Master cho kì thu tháng 9/2023.
Thông tin như sau:
1/ Yêu cầu: có kiến thức nền tảng về 1 trong các mảng như Edge Computing, Container Orchestration và Kubernetes, Machine Learning, Federated Learning và Edge AI. Tiếng anh yêu cầu có thể giao tiếp được với giáo sư trong quá trình phỏng vấn.
2/ Quyền lợi (Học phí~3,000,000 won/ 1 kì => tương ứng mỗi tháng 500,000 won):
👉Master: 1,400,000 won/ tháng
👉Phd:
+ 1,800,000 won/ tháng (2 năm đầu)
+ 2,000,000 won/ tháng (1 năm cuối)
👉Combined:
+ 1,600,000 won/ tháng (2 năm đầu)
+ 2,000,000 won/ tháng (3 năm cuối)
Ngoài ra cuối mỗi học kì giáo sư đều thưởng riêng nếu có performance tốt và publish paper.
3/ Thông tin khác:
+ Nộp CV, bảng điểm, research statement cho giáo sư Taehong Kim: taehongkim@cbnu.ac.kr (Deadline: 24/04/2023)
+ Lab info: https://nclab.cbnu.ac.kr
+ Guidelines: https://cia.chungbuk.ac.kr/.../notice/view/wr_id/138/key/180
Cảm ơn mọi người đã dành thời gian đọc tin 💖💖💖
=====================================
Network Computing Laboratory (NCLab) at CBNU is recruiting highly motivated PhD, MS-PhD Integrated, and MS students.
Requirements:
👉 A strong academic and research background in computer science, information technology, or a closely related discipline.
👉Keen research interests in one or more of the following areas: Edge Computing, Container Orchestration and Kubernetes, Machine Learning, Federated Learning, and Edge AI.
👉Good Python and Java programming abilities, as well as a firm grasp of Linux systems.
👉Experience in using tools like Git, GitHub, and open-source projects related to the research interests.
👉Must be able to independently and collaboratively propose, demonstrate, and evaluate new ideas.
Salary:
👉Master: 1,400,000 KRW/month
👉Ph.D: 1,800,000 KRW/month (for first 2 years), 2,000,000 KRW/month (from 3rd year)
👉Combined: 1,600,000 KRW/month (for first 2 years), same as Ph.D (from 3rd year)
👉Aside from the basic salary, incentives will be given based on achievement and contribution to the lab, as well as any travel expenses for presenting work at an international conference will also be provided.
Schedule:
👉Deadline of application submission: 21 April, 2023
👉Evaluation and interview, if needed: 22 April – 4 May, 2023
👉Notification of acceptance: 5 May, 2023
👉Official admission process: 15 – 26 May, 2023
* The applicants can see the detail guideline and the required application materials at https://cia.chungbuk.ac.kr/.../notice/view/wr_id/138/key/180. How to apply: If you are interested in our lab, please send your resume, transcripts, and a short research statement (about 300 to 500 words in total) to Prof. Taehong Kim (taehongkim@cbnu.ac.kr). For more details, you can visit our homepage at http://nclab.cbnu.ac.kr
'@
$ws.Range("E18").Value = $longText

# --- Row heights: rows 18-22 get a custom height of 185; row 23 keeps the default (set last so autofit does not override it) ---
$ws.Rows.Item(18).RowHeight = 185
$ws.Rows.Item(19).RowHeight = 185
$ws.Rows.Item(20).RowHeight = 185
$ws.Rows.Item(21).RowHeight = 185
$ws.Rows.Item(22).RowHeight = 185

# --- View state: scroll/zoom/selection to match the saved workbook view ---
$excel.ActiveWindow.Zoom = 68
$ws.Range("H18:H23").Select()

